$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width (new column added for Label text).
# Target OOXML width is 20.85546875 chars; the engine quantizes ColumnWidth
# in 1/12-character steps, so 20 lands on the closest achievable bucket
# (rendered width 20.8333... chars).
$ws.Columns.Item(2).ColumnWidth = 20

# Establish the shared-string insertion order to match the target workbook:
# Label, others, Test Case#, TC01, TC02, TC03, TC04, then the Trust* labels.
$ws.Range("B1").Value = "Label"
$ws.Range("G2").Value = "others"

$ws.Range("A1").Value = "Test Case#"
$ws.Range("A2").Value = "TC01"
$ws.Range("A3").Value = "TC02"
$ws.Range("A4").Value = "TC03"
$ws.Range("A5").Value = "TC04"

$ws.Range("B3").Value = "Trust RR8 vs Region R1"
$ws.Range("B4").Value = "Trust RR8 vs Peers"
$ws.Range("B5").Value = "Trust RR1 vs Peers"

# Remaining, already-existing shared strings reused as-is
$ws.Range("C1").Value = "Chart3"
$ws.Range("E1").Value = "UserName"
$ws.Range("F1").Value = "Password"
$ws.Range("G1").Value = "Label"
$ws.Range("E2").Value = "admin"
$ws.Range("F2").Value = "admin"

$ws.Range("B2").Value = ""
$ws.Range("C5").Clear()

# Apply the same border style (style index 1) to the new cells in column A and B
$ws.Range("A1:B5").Borders.LineStyle = 1

$ws.Range("B5").Select()
